$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to keep a numeric-looking literal as TEXT (matches the
    # original inlineStr cell), then restore the default "Normal" style so
    # no stray NumberFormat/style is left behind on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" "70.646.86"
$ws.Range("E2").Value = "  +2.55%  "
Set-TextValue $ws "D3" "3.572.97"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws "D5" "597.98"
$ws.Range("E5").Value = "  +1.76%  "
Set-TextValue $ws "D6" "173.09"
$ws.Range("E6").Value = "  +1.90%  "
Set-TextValue $ws "D7" "3.567.34"
$ws.Range("E7").Value = "  +1.71%  "
Set-TextValue $ws "D8" "0.616"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +6.33%  "
Set-TextValue $ws "D11" "7.41"
$ws.Range("E11").Value = "  +9.31%  "
Set-TextValue $ws "D12" "0.590"
$ws.Range("E12").Value = "  +2.29%  "
Set-TextValue $ws "D13" "46.66"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("E14").Value = "  +1.53%  "
Set-TextValue $ws "D15" "4.149.30"
$ws.Range("E15").Value = "  +1.62%  "
Set-TextValue $ws "D16" "8.40"
$ws.Range("E16").Value = "  -0.01%  "
Set-TextValue $ws "D17" "613.39"
$ws.Range("E17").Value = "  +0.46%  "
Set-TextValue $ws "D18" "3.558.59"
$ws.Range("E18").Value = "  +1.30%  "
Set-TextValue $ws "D19" "70.667.87"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("E20").Value = "  -0.90%  "
Set-TextValue $ws "D21" "17.43"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("E22").Value = "  +0.59%  "
Set-TextValue $ws "D23" "9.24"
$ws.Range("E23").Value = "  -16.88%  "
Set-TextValue $ws "D24" "15.80"
$ws.Range("E24").Value = "  +1.03%  "
Set-TextValue $ws "D25" "97.03"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("E27").Value = "  +0.07%  "
Set-TextValue $ws "D28" "2.64"
$ws.Range("E28").Value = "  +1.26%  "
Set-TextValue $ws "D29" "33.90"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("E30").Value = "  -0.06%  "
Set-TextValue $ws "D31" "8.36"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  -1.76%  "
Set-TextValue $ws "D33" "7.18"
$ws.Range("E33").Value = "  +4.42%  "
Set-TextValue $ws "D34" "661.02"
$ws.Range("E34").Value = "  +7.04%  "
$ws.Range("E35").Value = "  -0.91%  "
Set-TextValue $ws "D36" "3.69"
$ws.Range("E36").Value = "  +7.09%  "
$ws.Range("E37").Value = "  -0.66%  "
Set-TextValue $ws "D38" "10.83"
$ws.Range("E38").Value = "  +1.22%  "
Set-TextValue $ws "D39" "0.0478"
$ws.Range("E39").Value = "  +8.26%  "
Set-TextValue $ws "D40" "57.35"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +5.81%  "
Set-TextValue $ws "D43" "3.392.50"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  -0.69%  "
Set-TextValue $ws "D45" "0.0₃0714"
$ws.Range("E45").Value = "  +3.32%  "
Set-TextValue $ws "D46" "32.93"
$ws.Range("E46").Value = "  +1.14%  "
Set-TextValue $ws "D47" "2.95"
$ws.Range("E47").Value = "  +7.38%  "
Set-TextValue $ws "D48" "2.66"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("E49").Value = "  +1.00%  "
Set-TextValue $ws "D50" "132.34"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  -0.06%  "
